$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix markdown rendering: cell text used the literal "<br>" sequence as a
# line-break marker (meant for markdown), which never actually renders as
# a line break. Replace each "<br>" with a genuine newline (preceded by a
# space, matching the source edit) so Excel shows real line breaks.
[void]$ws.Cells.Replace("<br>", " `n")

# Replacing the text causes the affected rows to auto-grow; re-run AutoFit
# so the row heights settle back to their original (default) sizing
# instead of staying pinned to an explicit custom height.
$ws.Range("A2:A37").EntireRow.AutoFit()

# Refreshed underlying statistics for a few rows (mean / CI_low / CI_high).
$ws.Range("B2").Value = 73.6593446835017
$ws.Range("C2").Value = 72.2911397665722
$ws.Range("D2").Value = 75.0275496004312

$ws.Range("B12").Value = 77.4373747871973
$ws.Range("C12").Value = 73.0041637008231
$ws.Range("D12").Value = 81.8705858735715

$ws.Range("B14").Value = 68.7195237359994
$ws.Range("C14").Value = 67.2821385114178
$ws.Range("D14").Value = 70.1569089605811

$ws.Range("B24").Value = 69.7923192501879
$ws.Range("C24").Value = 64.9239512310953
$ws.Range("D24").Value = 74.6606872692806
